$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New A-column model names (reordered) for rows 2..26
$names = @("model_34_9_0", "model_34_9_22", "model_34_9_21", "model_34_9_20", "model_34_9_19", "model_34_9_18", "model_34_9_17", "model_34_9_16", "model_34_9_15", "model_34_9_14", "model_34_9_13", "model_34_9_23", "model_34_9_12", "model_34_9_10", "model_34_9_9", "model_34_9_8", "model_34_9_7", "model_34_9_6", "model_34_9_5", "model_34_9_4", "model_34_9_3", "model_34_9_2", "model_34_9_1", "model_34_9_11", "model_34_9_24")

# Constant metric values (columns B..Q) shared by every row after the edit
$Bval = 0.9997058837929482
$Cval = 0.7014553523661253
$Dval = 0.9997041626281999
$Eval = 0.9999888437131497
$Fval = 0.9999250027496585
$Gval = 0.0001746000856645352
$Hval = 0.1772289992927379
$Ival = 0.0001424082685028334
$Jval = 0.00000632196281124164
$Kval = 0.00007436723823477589
$Lval = 0.005296279504786582
$Mval = 0.01321363256884855
$Nval = 1.000243406516181
$Oval = 0.01377616426329468
$Pval = 123.3060248478516
$Qval = 187.9064435658663

for ($i = 2; $i -le 26; $i++) {
    $ws.Cells.Item($i, 1).Value = $names[$i - 2]
    $ws.Cells.Item($i, 2).Value = $Bval
    $ws.Cells.Item($i, 3).Value = $Cval
    $ws.Cells.Item($i, 4).Value = $Dval
    $ws.Cells.Item($i, 5).Value = $Eval
    $ws.Cells.Item($i, 6).Value = $Fval
    $ws.Cells.Item($i, 7).Value = $Gval
    $ws.Cells.Item($i, 8).Value = $Hval
    $ws.Cells.Item($i, 9).Value = $Ival
    $ws.Cells.Item($i, 10).Value = $Jval
    $ws.Cells.Item($i, 11).Value = $Kval
    $ws.Cells.Item($i, 12).Value = $Lval
    $ws.Cells.Item($i, 13).Value = $Mval
    $ws.Cells.Item($i, 14).Value = $Nval
    $ws.Cells.Item($i, 15).Value = $Oval
    $ws.Cells.Item($i, 16).Value = $Pval
    $ws.Cells.Item($i, 17).Value = $Qval
}

